$d = $word.ActiveDocument

$replacements = @(
    @{old="946×8=7568"; new="517×2=1034"},
    @{old="845×4=3380"; new="708×9=6372"},
    @{old="706×9=6354"; new="620×9=5580"},
    @{old="955×5=4775"; new="827×9=7443"},
    @{old="279×8=2232"; new="772×2=1544"},
    @{old="253×8=2024"; new="218×9=1962"},
    @{old="367×4=1468"; new="421×4=1684"},
    @{old="323×5=1615"; new="358×7=2506"},
    @{old="935×8=7480"; new="187×9=1683"},
    @{old="801×9=7209"; new="956×5=4780"},
    @{old="936×5=4680"; new="421×3=1263"},
    @{old="692×3=2076"; new="573×7=4011"},
    @{old="892×7=6244"; new="745×9=6705"},
    @{old="241×8=1928"; new="694×7=4858"},
    @{old="977×4=3908"; new="318×5=1590"},
    @{old="446×8=3568"; new="388×5=1940"},
    @{old="217×2=434";  new="307×6=1842"},
    @{old="332×3=996";  new="132×7=924"},
    @{old="370×8=2960"; new="367×6=2202"},
    @{old="376×4=1504"; new="455×9=4095"},
    @{old="422×9=3798"; new="303×9=2727"},
    @{old="931×6=5586"; new="886×9=7974"},
    @{old="150×9=1350"; new="593×5=2965"},
    @{old="493×3=1479"; new="890×4=3560"},
    @{old="839×7=5873"; new="688×5=3440"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
